$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for column G, rows 2-36
$newValues = @{
    2  = 3
    3  = 3
    4  = 7
    5  = 5
    6  = 9
    7  = 7
    8  = 6
    9  = 6
    10 = 3
    11 = 2
    12 = 7
    13 = 5
    14 = 4
    15 = 6
    16 = 4
    17 = 4
    18 = 6
    19 = 8
    20 = 1
    21 = 4
    22 = 5
    23 = 3
    24 = 5
    25 = 2
    26 = 1
    27 = 6
    28 = 4
    29 = 7
    30 = 7
    31 = 2
    32 = 5
    33 = 9
    34 = 3
    35 = 4
    36 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
